$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header change
$ws.Range("A1").Value = "Pattern"

# Full data block (rows 2-22), columns A:E
$data = @(
    @("ATATATATAT", 0, 0.667, 0, 0.22),
    @("ATATATAAAT", 1, 1.667, 1, 1.22),
    @("AAATATATAT", 1, 1.667, 1, 1.22),
    @("ATATAAATAT", 1, 2, 1, 1.33),
    @("TAATATATAT", 1, 3.667, 0, 1.56),
    @("ATAAATATAT", 1, 3, 1, 1.67),
    @("AAATATAAAT", 2, 2.667, 2, 2.22),
    @("ATATAAAAAT", 2, 3, 2, 2.33),
    @("AAATAAATAT", 2, 3, 2, 2.33),
    @("ATATAATAAT", 2, 5, 0, 2.33),
    @("TAATATAAAT", 2, 4.667, 1, 2.56),
    @("TAATAAATAT", 2, 5, 1, 2.67),
    @("AAAAATATAT", 2, 4, 2, 2.67),
    @("ATAAATAAAT", 2, 4, 2, 2.67),
    @("AATAATATAT", 2, 7, 0, 3),
    @("TAAAATATAT", 2, 6, 2, 3.33),
    @("AAATAAAAAT", 3, 4, 3, 3.33),
    @("AAATAATAAT", 3, 6, 1, 3.33),
    @("TAATAAAAAT", 3, 6, 2, 3.67),
    @("AAAAATAAAT", 3, 5, 3, 3.67),
    @("TAAAATAAAT", 3, 7, 3, 4.33)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $row++
}
